$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-03 10:09:45"

# ---------------------------------------------------------------------------
# Sheet "Главные"
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Главные")

# Refresh the as_of_utc timestamp for every data row (2-26)
$wsMain.Range("AA2:AA26").Value = $newTimestamp

# Row 8 - Gamaley Evgeniy
$wsMain.Range("C8").Value = 28
$wsMain.Range("D8").Value = 510
$wsMain.Range("E8").Value = 271
$wsMain.Range("F8").Value = 239
$wsMain.Range("G8").Value = 18.21
$wsMain.Range("H8").Value = 9.68
$wsMain.Range("I8").Value = 8.539999999999999
$wsMain.Range("J8").Value = 118
$wsMain.Range("K8").Value = 112
$wsMain.Range("V8").Value = 14

# Row 9 - Gashilov Viktor
$wsMain.Range("C9").Value = 31
$wsMain.Range("D9").Value = 488
$wsMain.Range("E9").Value = 255
$wsMain.Range("F9").Value = 233
$wsMain.Range("G9").Value = 15.74
$wsMain.Range("H9").Value = 8.23
$wsMain.Range("I9").Value = 7.52
$wsMain.Range("J9").Value = 125
$wsMain.Range("K9").Value = 114
$wsMain.Range("V9").Value = 20

# Row 11 - Dudarov Aleksandr
$wsMain.Range("C11").Value = 23
$wsMain.Range("D11").Value = 548
$wsMain.Range("E11").Value = 252
$wsMain.Range("F11").Value = 296
$wsMain.Range("G11").Value = 23.83
$wsMain.Range("H11").Value = 10.96
$wsMain.Range("I11").Value = 12.87
$wsMain.Range("J11").Value = 111
$wsMain.Range("K11").Value = 103
$wsMain.Range("U11").Value = 1
$wsMain.Range("V11").Value = 16

# Row 18 - Naumov Denis
$wsMain.Range("C18").Value = 29
$wsMain.Range("D18").Value = 464
$wsMain.Range("E18").Value = 233
$wsMain.Range("F18").Value = 231
$wsMain.Range("G18").Value = 16
$wsMain.Range("H18").Value = 8.029999999999999
$wsMain.Range("I18").Value = 7.97
$wsMain.Range("J18").Value = 94
$wsMain.Range("K18").Value = 108
$wsMain.Range("V18").Value = 8

# Row 19 - Ovchinnikov Pavel
$wsMain.Range("C19").Value = 24
$wsMain.Range("D19").Value = 416
$wsMain.Range("E19").Value = 208
$wsMain.Range("F19").Value = 208
$wsMain.Range("G19").Value = 17.33
$wsMain.Range("H19").Value = 8.67
$wsMain.Range("I19").Value = 8.67
$wsMain.Range("J19").Value = 99
$wsMain.Range("K19").Value = 89
$wsMain.Range("V19").Value = 12

# Row 21 - Romasko Evgeniy
$wsMain.Range("C21").Value = 27
$wsMain.Range("D21").Value = 378
$wsMain.Range("E21").Value = 168
$wsMain.Range("F21").Value = 210
$wsMain.Range("G21").Value = 14
$wsMain.Range("H21").Value = 6.22
$wsMain.Range("I21").Value = 7.78
$wsMain.Range("J21").Value = 74
$wsMain.Range("K21").Value = 90
$wsMain.Range("U21").Value = 2
$wsMain.Range("V21").Value = 8

# ---------------------------------------------------------------------------
# Sheet "Линейные"
# ---------------------------------------------------------------------------
$wsLinear = $wb.Worksheets.Item("Линейные")

# Refresh the as_of_utc timestamp for every data row (2-26)
$wsLinear.Range("AA2:AA26").Value = $newTimestamp

# Row 3 - Bersenyov Maksim
$wsLinear.Range("C3").Value = 30
$wsLinear.Range("D3").Value = 457
$wsLinear.Range("E3").Value = 249
$wsLinear.Range("F3").Value = 208
$wsLinear.Range("G3").Value = 15.23
$wsLinear.Range("H3").Value = 8.300000000000001
$wsLinear.Range("I3").Value = 6.93
$wsLinear.Range("J3").Value = 112
$wsLinear.Range("K3").Value = 89
$wsLinear.Range("V3").Value = 14

# Row 12 - Zaytsev Valentin
$wsLinear.Range("C12").Value = 26
$wsLinear.Range("D12").Value = 454
$wsLinear.Range("E12").Value = 220
$wsLinear.Range("F12").Value = 234
$wsLinear.Range("G12").Value = 17.46
$wsLinear.Range("H12").Value = 8.460000000000001
$wsLinear.Range("I12").Value = 9
$wsLinear.Range("J12").Value = 100
$wsLinear.Range("K12").Value = 107
$wsLinear.Range("U12").Value = 2
$wsLinear.Range("V12").Value = 18

# Row 26 - Slavikovskiy Roman
$wsLinear.Range("C26").Value = 28
$wsLinear.Range("D26").Value = 566
$wsLinear.Range("E26").Value = 244
$wsLinear.Range("F26").Value = 322
$wsLinear.Range("G26").Value = 20.21
$wsLinear.Range("H26").Value = 8.710000000000001
$wsLinear.Range("I26").Value = 11.5
$wsLinear.Range("J26").Value = 97
$wsLinear.Range("K26").Value = 101
$wsLinear.Range("V26").Value = 8
